# Applies updated "想去人数" (F column) figures across sheets, per commit
# "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 152
$ws.Range("F5").Value = 2996
$ws.Range("F6").Value = 803
$ws.Range("F7").Value = 597
$ws.Range("F9").Value = 445
$ws.Range("F12").Value = 551
$ws.Range("F17").Value = 25
$ws.Range("F19").Value = 2682
$ws.Range("F25").Value = 648
$ws.Range("F29").Value = 13
$ws.Range("F33").Value = 126
$ws.Range("F34").Value = 905
$ws.Range("F35").Value = 4697
$ws.Range("F36").Value = 258
$ws.Range("F37").Value = 36
$ws.Range("F38").Value = 8
$ws.Range("F39").Value = 83

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 4
$ws.Range("F9").Value = 356
$ws.Range("F14").Value = 169
$ws.Range("F23").Value = 276
$ws.Range("F24").Value = 23
$ws.Range("F25").Value = 304
$ws.Range("F27").Value = 183
$ws.Range("F37").Value = 561

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1466
$ws.Range("F6").Value = 267
$ws.Range("F7").Value = 266

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1466
$ws.Range("F5").Value = 152
$ws.Range("F6").Value = 267
$ws.Range("F9").Value = 2996
$ws.Range("F10").Value = 803
$ws.Range("F11").Value = 597
$ws.Range("F13").Value = 445
$ws.Range("F16").Value = 551
$ws.Range("F18").Value = 356
$ws.Range("F25").Value = 169
$ws.Range("F27").Value = 2682
$ws.Range("F35").Value = 266
$ws.Range("F37").Value = 648
$ws.Range("F38").Value = 648
$ws.Range("F39").Value = 276
$ws.Range("F42").Value = 304
$ws.Range("F45").Value = 905
$ws.Range("F47").Value = 4697
$ws.Range("F48").Value = 258
$ws.Range("F50").Value = 561
